$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$Cvals = @(0.0464788531968594, 0.04138572481257086, 0.03828215826391101, 0.03702314928105466, 0.03681443078552604, 0.03826515596347235, 0.04471769998190211, 0.05756949035897208, 0.06714874068239851, 0.07153996008926811, 0.07320791032182683, 0.07284845672596418, 0.07167707980194393, 0.07096024849910521, 0.06686246206585622, 0.06435738649572897, 0.06291967697424639, 0.06243342524111029, 0.06462372881780709, 0.07202100130933786, 0.07688534585773255, 0.0742863383788972, 0.06450330771252766, 0.05406985807610454)
$Dvals = @(0.01438901658912073, 0.01373699565459674, 0.01334850134078636, 0.0131930362478343, 0.01316738970684028, 0.01334639331802734, 0.01416166171317457, 0.01586037637054005, 0.01717775731424354, 0.01779396702136893, 0.01802989508886554, 0.01797896649105013, 0.01781332434070038, 0.01771220455936628, 0.01713783983142747, 0.01678992608673013, 0.01659139946440291, 0.01652444934746455, 0.01682679707397483, 0.01786190616024186, 0.01855353528830506, 0.01818296601341984, 0.0168101230598694, 0.01538926631313586)
$Evals = @(0.0882776211624261, 0.08799679838679708, 0.08789105724109803, 0.08786472243632204, 0.08786136122106925, 0.08789063425299659, 0.08816694407802572, 0.0892387982657219, 0.09035109224365456, 0.09092806776114415, 0.09115679531851484, 0.09110707881148983, 0.09094667991023186, 0.09084976541531731, 0.09031481661710572, 0.09000484598890779, 0.08983323869378523, 0.08977628172631924, 0.09003715135753509, 0.09099351479439832, 0.0916782536179781, 0.09130732189785817, 0.09002252554843793, 0.08889192644198829)
$Fvals = @(5.642909543315625, 5.482549938652738, 5.386418769756347, 5.347819296562534, 5.341444272237879, 5.385895891294979, 5.587127066978752, 6.000762038580348, 6.317029612363683, 6.463766043766498, 6.51975691826334, 6.507679170748418, 6.468363854121037, 6.444337822330567, 6.307498795200559, 6.224295402088899, 6.176707858859317, 6.160641316244693, 6.233124627007442, 6.479900076935792, 6.643665388925115, 6.556029169001533, 6.229132165444582, 5.886747220328772)
$Gvals = @(0.00261952693826909, 0.00262730893706258, 0.002632327315738799, 0.002634433001328191, 0.002634786319817851, 0.002632355467659608, 0.002622160472378843, 0.002604062284101317, 0.002591903802676526, 0.002586616205130655, 0.002584648647593976, 0.002585070855119829, 0.002586453638210584, 0.00258730514925652, 0.002592254235163279, 0.002595352488051548, 0.00259715744223829, 0.002597772513241909, 0.002595020303330963, 0.002586046540417481, 0.002580384046207985, 0.002583387792103127, 0.002595170410030454, 0.002608757241864421)
$Jvals = @(0.2621564063342419, 0.255853961227956, 0.252137334713467, 0.2506608812822577, 0.2504180071465356, 0.2521172689899487, 0.2599513042672612, 0.2765477953987272, 0.2895219999239345, 0.2956002525491215, 0.2979277763925268, 0.2974253462851095, 0.2957912190368006, 0.294793646115906, 0.2891283586996849, 0.2856983614980777, 0.2837420751835538, 0.2830825395342771, 0.2860617734862956, 0.2962704971407675, 0.3030932275195397, 0.2994378531698203, 0.2858974261703509, 0.2719229922770268)
$Kvals = @(3.437531336783763, 3.305625578365778, 3.227675352603285, 3.196668959954195, 3.19156605141302, 3.227254123792306, 3.391415718970961, 3.737743048217624, 4.007504304577481, 4.133651315647967, 4.18192115737736, 4.171502995240076, 4.137612426947157, 4.11691890308191, 3.999329964085462, 3.928076894658034, 3.887416528101767, 3.873704800771691, 3.935628474525743, 4.147553259886649, 4.288980434205826, 4.213228239641694, 3.932213459010029, 3.641395263629818)

for ($i = 0; $i -lt 24; $i++) {
    $row = $i + 2
    $ws.Range("C" + $row).Value = [double]$Cvals[$i]
    $ws.Range("D" + $row).Value = [double]$Dvals[$i]
    $ws.Range("E" + $row).Value = [double]$Evals[$i]
    $ws.Range("F" + $row).Value = [double]$Fvals[$i]
    $ws.Range("G" + $row).Value = [double]$Gvals[$i]
    $ws.Range("J" + $row).Value = [double]$Jvals[$i]
    $ws.Range("K" + $row).Value = [double]$Kvals[$i]
}
